# GreyBoxTesting.xlsx edit script
# - Rename "LeaguesController" -> "LeagueController"
# - Append 20 new test case rows (TC_GBT_10 .. TC_GBT_29) for Match/Player/Stats/Team controllers
# - Update sheet selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the "LeaguesController" typo -> "LeagueController" for every row
#    in the League Controller block (rows 5-12, column C)
# ---------------------------------------------------------------------------
$ws.Range("C5:C12").Value = "LeagueController"

# ---------------------------------------------------------------------------
# 2) Append the new rows (13-32). Each row re-uses the same look & feel
#    (borders / alignment / wrap / row height) as the existing test-case
#    rows above it, so we copy formatting from a matching template row and
#    then overwrite the cell values.
# ---------------------------------------------------------------------------

# Row number -> template row that has the same column-A alignment style
$templateRow = @{
    13 = 5;  14 = 6;  15 = 7;  16 = 5;  17 = 5;
    18 = 5;  19 = 6;  20 = 7;  21 = 5;  22 = 5;
    23 = 5;  24 = 6;  25 = 7;  26 = 5;  27 = 5;
    28 = 5;  29 = 6;  30 = 7;  31 = 5;  32 = 5;
}

$rowValues = @{
    13 = @("Index Page in Match Controller returns a valid response",  "TC_GBT_10", "MatchController",  "Database Context", "Valid/Non Null Response", "Y")
    14 = @("Details Page in Match Controller returns a valid response", "TC_GBT_11", "MatchController",  "Database Context", "Valid/Non Null Response", "Y")
    15 = @("Create Page in Match Controller returns a valid response",  "TC_GBT_12", "MatchController",  "Database Context", "Valid/Non Null Response", "Y")
    16 = @("Edit Page in Match Controller returns a valid response",    "TC_GBT_13", "MatchController",  "Database Context", "Valid/Non Null Response", "Y")
    17 = @("Delete Page in Match Controller returns a valid response",  "TC_GBT_14", "MatchController",  "Database Context", "Valid/Non Null Response", "Y")
    18 = @("Index Page in Player Controller returns a valid response", "TC_GBT_15", "PlayerController", "Database Context", "Valid/Non Null Response", "Y")
    19 = @("Details Page in Player Controller returns a valid response","TC_GBT_16", "PlayerController", "Database Context", "Valid/Non Null Response", "Y")
    20 = @("Create Page in Player Controller returns a valid response", "TC_GBT_17", "PlayerController", "Database Context", "Valid/Non Null Response", "Y")
    21 = @("Edit Page in Player Controller returns a valid response",   "TC_GBT_18", "PlayerController", "Database Context", "Valid/Non Null Response", "Y")
    22 = @("Delete Page in Player Controller returns a valid response", "TC_GBT_19", "PlayerController", "Database Context", "Valid/Non Null Response", "Y")
    23 = @("Index Page in Stats Controller returns a valid response",  "TC_GBT_20", "StatsController",  "Database Context", "Valid/Non Null Response", "Y")
    24 = @("Details Page in Stats Controller returns a valid response","TC_GBT_21", "StatsController",  "Database Context", "Valid/Non Null Response", "Y")
    25 = @("Create Page in Stats Controller returns a valid response", "TC_GBT_22", "StatsController",  "Database Context", "Valid/Non Null Response", "Y")
    26 = @("Edit Page in Stats Controller returns a valid response",   "TC_GBT_23", "StatsController",  "Database Context", "Valid/Non Null Response", "Y")
    27 = @("Delete Page in Stats Controller returns a valid response", "TC_GBT_24", "StatsController",  "Database Context", "Valid/Non Null Response", "Y")
    28 = @("Index Page in Team Controller returns a valid response",   "TC_GBT_25", "TeamController",   "Database Context", "Valid/Non Null Response", "Y")
    29 = @("Details Page in Team Controller returns a valid response", "TC_GBT_26", "TeamController",   "Database Context", "Valid/Non Null Response", "Y")
    30 = @("Create Page in Team Controller returns a valid response",  "TC_GBT_27", "TeamController",   "Database Context", "Valid/Non Null Response", "Y")
    31 = @("Edit Page in Team Controller returns a valid response",    "TC_GBT_28", "TeamController",   "Database Context", "Valid/Non Null Response", "Y")
    32 = @("Delete Page in Team Controller returns a valid response",  "TC_GBT_29", "TeamController",   "Database Context", "Valid/Non Null Response", "Y")
}

for ($r = 13; $r -le 32; $r++) {
    $src = $templateRow[$r]

    # Copy the whole template row's formatting into the new row
    $ws.Range("A" + $src + ":F" + $src).Copy()
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122)

    # Make sure row height + wrap text match the rest of the test-case table
    $ws.Rows.Item($r).RowHeight = 57.45
    $ws.Range("A" + $r + ":F" + $r).WrapText = $true

    $vals = $rowValues[$r]
    $ws.Range("A" + $r).Value = $vals[0]
    $ws.Range("B" + $r).Value = $vals[1]
    $ws.Range("C" + $r).Value = $vals[2]
    $ws.Range("D" + $r).Value = $vals[3]
    $ws.Range("E" + $r).Value = $vals[4]
    $ws.Range("F" + $r).Value = $vals[5]
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Update the view so the selection / scroll position matches the final
#    state of the sheet (top-left back to A1, active cell C31)
# ---------------------------------------------------------------------------
[void]$ws.Range("A1").Select()
[void]$ws.Range("C31").Select()
